# B6-PowerPoint.pptx edit: re-style the three small tables (slides 14-16)
# from the bespoke "Table_0" style to the built-in "Medium Style 2 - Accent 1"
# table style ({24DE159C-679A-4F66-BA1B-85144E7A756C}), matching what the
# Table Design gallery writes into <a:tableStyleId> for each <a:tbl>.

$p = $ppt.ActivePresentation

$newTableStyleId = "{24DE159C-679A-4F66-BA1B-85144E7A756C}"

# Slides 14, 15 and 16 (1-based, matching Slides.Item order) each hold a
# table as their very first shape (a graphicFrame placed before any other
# shape in the slide's shape tree).
$slideIndexesWithTables = 14, 15, 16

foreach ($slideIndex in $slideIndexesWithTables) {
    $slide = $p.Slides.Item($slideIndex)
    $tableShape = $slide.Shapes.Item(1)

    if ($tableShape.HasTable) {
        $table = $tableShape.Table
        $table.ApplyStyle($newTableStyleId)
    }
}
